$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35-113 down to 36-114.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly data point.
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value2 = 45028
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = 100112012
$ws.Range("G35").Value = "Espinaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 8000
$ws.Range("L35").Value = 9000
$ws.Range("M35").Value = 8375
$ws.Range("N35").Value = "$/cuna 10 kilos"
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 838
$ws.Range("Q35").Value = 10
$ws.Range("R35").Value = "Hortaliza"
